$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. screenTitles (sheet4): append three new rows (MyPreference / Menu / home)
# ---------------------------------------------------------------------------
$wsTitles = $wb.Worksheets.Item("screenTitles")

$wsTitles.Range("A11").Value = "MyPreference"
$wsTitles.Range("A11").Font.Name = "Calibri"
$wsTitles.Range("B11").Value = "mijn voorkeur"
$wsTitles.Range("B11").Font.Color = 0x222222

$wsTitles.Range("A12").Value = "Menu"
$wsTitles.Range("A12").Font.Name = "Calibri"
$wsTitles.Range("B12").Value = "taal van het menu"
$wsTitles.Range("B12").Font.Color = 0x222222

$wsTitles.Range("A13").Value = "home"
$wsTitles.Range("A13").Font.Name = "Calibri"
$wsTitles.Range("B13").Value = "home"
$wsTitles.Range("B13").Font.Color = 0x222222

# ---------------------------------------------------------------------------
# 2. Add the new "parameters" worksheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsParams = $wb.Worksheets.Add($null, $lastSheet)
$wsParams.Name = "parameters"

$wsParams.Columns.Item(1).ColumnWidth = 28.7369791666667
$wsParams.Columns.Item(2).ColumnWidth = 32.7369791666667
$wsParams.Columns.Item(3).ColumnWidth = 16.8776041666667

# Header row (re-uses existing shared strings "objectID" / "name_nl" and the
# bold/shaded header style already used on the other lookup sheets)
$wsParams.Range("A1").Value = "objectID"
$wsParams.Range("A1").Font.Bold = $true
$wsParams.Range("A1").Interior.Color = 0xD9D9D9
$wsParams.Range("B1").Value = "name_nl"
$wsParams.Range("B1").Font.Bold = $true
$wsParams.Range("B1").Interior.Color = 0xD9D9D9

$wsParams.Range("A2").Value = "language_FR"
$wsParams.Range("A2").Font.Color = 0x222222
$wsParams.Range("B2").Value = "FR"
$wsParams.Range("B2").Font.Name = "Calibri"

$wsParams.Range("A3").Value = "langauge_NL"
$wsParams.Range("A3").Font.Color = 0x222222
$wsParams.Range("B3").Value = "NL"
$wsParams.Range("B3").Font.Name = "Calibri"

$wsParams.Range("A4").Value = "preferenceNL"
$wsParams.Range("A4").Font.Name = "Calibri"
$wsParams.Range("B4").Value = "mijn voorkeur"
$wsParams.Range("B4").Font.Name = "Calibri"

$wsParams.Range("A5").Value = "preferenceFrench"
$wsParams.Range("A5").Font.Name = "Calibri"
$wsParams.Range("B5").Value = "préférences"
$wsParams.Range("B5").Font.Name = "Calibri"

$wsParams.Range("A6").Value = "confirmNL"
$wsParams.Range("A6").Font.Name = "Calibri"
$wsParams.Range("B6").Value = "bevestigen"
$wsParams.Range("B6").Font.Color = 0x222222

$wsParams.Range("A7").Value = "confirmFR"
$wsParams.Range("A7").Font.Name = "Calibri"
$wsParams.Range("B7").Value = "confirmer"
$wsParams.Range("B7").Font.Color = 0x222222

$wsParams.Range("A8").Value = "languageMenuFR"
$wsParams.Range("A8").Font.Name = "Calibri"
$wsParams.Range("B8").Value = "langue du menu"
$wsParams.Range("B8").Font.Color = 0x222222

$wsParams.Range("A9").Value = "languageMenuNL"
$wsParams.Range("A9").Font.Name = "Calibri"
$wsParams.Range("B9").Value = "taal van het menu"
$wsParams.Range("B9").Font.Color = 0x222222

$wsParams.Range("A10").Value = "language_the_shop_FR"
$wsParams.Range("A10").Font.Color = 0x222222
$wsParams.Range("B10").Value = "langue du shop"
$wsParams.Range("B10").Font.Color = 0x222222

$wsParams.Range("A11").Value = "language_the_shop_NL"
$wsParams.Range("A11").Font.Color = 0x222222
$wsParams.Range("B11").Value = "taal van de shop"
$wsParams.Range("B11").Font.Color = 0x222222

# trailing blank styled rows
$wsParams.Range("A12").Font.Color = 0x222222
$wsParams.Range("B12").Font.Color = 0x222222
$wsParams.Range("A13").Font.Color = 0x222222
$wsParams.Range("B13").Font.Color = 0x222222
$wsParams.Range("A14").Font.Color = 0x222222
$wsParams.Range("B14").Font.Color = 0x222222

$wsParams.PageSetup.Orientation = 1

$wsParams.Range("B14").Select()

# ---------------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping so the saved view matches
# ---------------------------------------------------------------------------
$wsFilms = $wb.Worksheets.Item("Films")
$wsFilms.Range("B8").Select()

$wsTitles.Activate()
$wsTitles.Range("H16").Select()
